$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update selection on the first sheet (matches diff: tabSelected removed, topLeftCell removed, selection -> E1)
$ws1.Range("E1").Select()

# Add the new sheet after WD_person_matches
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$newSheet.Name = "previously_identified_matches"

$newSheet.Range("B1").Value = "84000 ID"
$newSheet.Range("C1").Value = "BDRC ID"
$newSheet.Range("B2").Value = "eft:sarvajnadeva"
$newSheet.Range("C2").Value = "P00KG07267"
$newSheet.Range("B3").Value = "eft:vidyakaraprabha"
$newSheet.Range("C3").Value = "P8211"
$newSheet.Range("B4").Value = "eft:dharmakara"
$newSheet.Range("C4").Value = "P8249"
$newSheet.Range("B5").Value = "eft:jinamitra-k-"
$newSheet.Range("C5").Value = "P8209"
$newSheet.Range("B6").Value = "eft:klu-i-rgyal-mtshan"
$newSheet.Range("C6").Value = "P8183"
$newSheet.Range("B7").Value = "eft:dzi-na-mi-tra-k-"
$newSheet.Range("C7").Value = "P8209"
$newSheet.Range("B8").Value = "eft:cog-ro-klu-i-rgyal-mtshan"
$newSheet.Range("C8").Value = "P8183"
$newSheet.Range("B9").Value = "eft:ban-de-dpal-gyi-lhun-po"
$newSheet.Range("C9").Value = "P4259"
$newSheet.Range("B10").Value = "eft:ban-de-dpal-brtsegs"
$newSheet.Range("C10").Value = "P8182"
$newSheet.Range("B11").Value = "eft:dpal-byor"
$newSheet.Range("C11").Value = "P4258"
$newSheet.Range("B12").Value = "eft:surendrabodhi"
$newSheet.Range("C12").Value = "P8228"
$newSheet.Range("B13").Value = "eft:jinamitra"
$newSheet.Range("C13").Value = "P8209"
$newSheet.Range("B14").Value = "eft:danasila"
$newSheet.Range("C14").Value = "P3214"
$newSheet.Range("B15").Value = "eft:munivarman"
$newSheet.Range("C15").Value = "P8261"
$newSheet.Range("B16").Value = "eft:prajnavarman"
$newSheet.Range("C16").Value = "P2548"
$newSheet.Range("B17").Value = "eft:dpal-dbyangs"
$newSheet.Range("C17").Value = "P8260"
$newSheet.Range("B18").Value = "eft:ska-ba-dpal-brtsegs"
$newSheet.Range("C18").Value = "P8182"
$newSheet.Range("B19").Value = "eft:silendrabodhi"
$newSheet.Range("C19").Value = "P1KG8854"
$newSheet.Range("B20").Value = "eft:prajnavarma"
$newSheet.Range("C20").Value = "P2548"
$newSheet.Range("B21").Value = "eft:dipamkarasrijnana"
$newSheet.Range("C21").Value = "P3379"
$newSheet.Range("B22").Value = "eft:buddhaprabha"
$newSheet.Range("C22").Value = "P8268"
$newSheet.Range("B23").Value = "eft:ye-shes-sde"
$newSheet.Range("C23").Value = "P8205"
$newSheet.Range("B24").Value = "eft:dgon-gling-rma"
$newSheet.Range("C24").Value = "P8269"
$newSheet.Range("B25").Value = "eft:dpal-gyi-lhun-po"
$newSheet.Range("C25").Value = "P4259"
$newSheet.Range("B26").Value = "eft:dpal-brtsegs"
$newSheet.Range("C26").Value = "P8182"
$newSheet.Range("B27").Value = "eft:sakyaprabha"
$newSheet.Range("C27").Value = "P4CZ16819"
$newSheet.Range("B28").Value = "eft:dharmatasila"
$newSheet.Range("C28").Value = "P8266"
$newSheet.Range("B29").Value = "eft:ye-shes-snying-po"
$newSheet.Range("C29").Value = "P4255"
$newSheet.Range("B30").Value = "eft:visuddhasimha"
$newSheet.Range("C30").Value = "P8219"
$newSheet.Range("B31").Value = "eft:dge-ba-dpal"
$newSheet.Range("C31").Value = "P4263"
$newSheet.Range("B32").Value = "eft:devacandra"
$newSheet.Range("C32").Value = "P8220"
$newSheet.Range("B33").Value = "eft:kamalagupta"
$newSheet.Range("C33").Value = "P8093"
$newSheet.Range("B34").Value = "eft:rin-chen-bzang-po"
$newSheet.Range("C34").Value = "P753"
$newSheet.Range("B35").Value = "eft:rin-chen-tsho"
$newSheet.Range("C35").Value = "P8273"
$newSheet.Range("B36").Value = "eft:jnanagarbha"
$newSheet.Range("C36").Value = "P8217"
$newSheet.Range("B37").Value = "eft:vijayasila"
$newSheet.Range("C37").Value = "P8267"
$newSheet.Range("B38").Value = "eft:hwa-shang-zab-mo"
$newSheet.Range("C38").Value = "P0TMP080"
$newSheet.Range("B39").Value = "eft:rnam-par-mi-rtog-pa"
$newSheet.Range("C39").Value = "P0TMPT007"
$newSheet.Range("B40").Value = "eft:munivarma"
$newSheet.Range("C40").Value = "P8261"
$newSheet.Range("B41").Value = "eft:ratnaraksita"
$newSheet.Range("C41").Value = "P8265"
$newSheet.Range("B42").Value = "eft:dharmasribhadra"
$newSheet.Range("C42").Value = "P8171"
$newSheet.Range("B43").Value = "eft:gayadhara"
$newSheet.Range("C43").Value = "P8151"
$newSheet.Range("B44").Value = "eft:krsnapandita"
$newSheet.Range("C44").Value = "P2956"
$newSheet.Range("B45").Value = "eft:tshul-khrims-rgyal-ba"
$newSheet.Range("C45").Value = "P3456"
$newSheet.Range("B46").Value = "eft:celu"
$newSheet.Range("C46").Value = "P8206"
$newSheet.Range("B47").Value = "eft:jnanasiddhi"
$newSheet.Range("C47").Value = "P8222"
$newSheet.Range("B48").Value = "eft:punyasambhava"
$newSheet.Range("C48").Value = "P0TMP104"
$newSheet.Range("B49").Value = "eft:palgyi-lh-npo"
$newSheet.Range("C49").Value = "P4259"
$newSheet.Range("B50").Value = "eft:sarvanyadeva"
$newSheet.Range("C50").Value = "P00KG07267"
$newSheet.Range("B51").Value = "eft:vidyakaraprabha"
$newSheet.Range("C51").Value = "P8211"
$newSheet.Range("B52").Value = "eft:jinamitra"
$newSheet.Range("C52").Value = "P8209"
$newSheet.Range("B53").Value = "eft:sakyasena"
$newSheet.Range("C53").Value = "?"
$newSheet.Range("B54").Value = "eft:jnanasiddhi"
$newSheet.Range("C54").Value = "P8222"
$newSheet.Range("B55").Value = "eft:anandasri-s-"
$newSheet.Range("C55").Value = "P0TMP092"
$newSheet.Range("B56").Value = "eft:kawa-paltsek-under-the-name-paltsek-raksita-"
$newSheet.Range("C56").Value = "P8182"
$newSheet.Range("B57").Value = "eft:yesh-d-"
$newSheet.Range("C57").Value = "P8205"
$newSheet.Range("B58").Value = "eft:danasila"
$newSheet.Range("C58").Value = "P3214"
$newSheet.Range("B59").Value = "eft:t-vidyakarasimha"
$newSheet.Range("C59").Value = "P8213"
$newSheet.Range("B60").Value = "eft:yesh-nyingpo"
$newSheet.Range("C60").Value = "P4255"
$newSheet.Range("B61").Value = "eft:band-yesh-d-"
$newSheet.Range("C61").Value = "P8205"
$newSheet.Range("B62").Value = "eft:leki-d-"
$newSheet.Range("C62").Value = "P8263"
$newSheet.Range("B63").Value = "eft:surendrabodhi"
$newSheet.Range("C63").Value = "P1KG8854"
$newSheet.Range("B64").Value = "eft:yesh-d-ye-shes-sde-"
$newSheet.Range("C64").Value = "P8205"
$newSheet.Range("B65").Value = "eft:srilendrabodhi"
$newSheet.Range("C65").Value = "P1KG8854"
$newSheet.Range("B66").Value = "eft:silendrabodhi"
$newSheet.Range("C66").Value = "P1KG8854"
$newSheet.Range("B67").Value = "eft:dipamkarasrijnana"
$newSheet.Range("C67").Value = "P3379"
$newSheet.Range("B68").Value = "eft:prajnavarma"
$newSheet.Range("C68").Value = "P2548"
$newSheet.Range("B69").Value = "eft:dharmapala"
$newSheet.Range("C69").Value = "P0RK8"
$newSheet.Range("B70").Value = "eft:t-jnanagarbha"
$newSheet.Range("C70").Value = "P4255"
$newSheet.Range("B71").Value = "eft:t-jnanagarbha"
$newSheet.Range("C71").Value = "P8217"
$newSheet.Range("B72").Value = "eft:vidyakarasimha"
$newSheet.Range("C72").Value = "P8213"
$newSheet.Range("B73").Value = "eft:dharmatasila"
$newSheet.Range("C73").Value = "P8266"
$newSheet.Range("B74").Value = "eft:ch-nyi-tsultrim"
$newSheet.Range("C74").Value = "P8266"
$newSheet.Range("B75").Value = "eft:jnanasidhi"
$newSheet.Range("C75").Value = "P8222"
$newSheet.Range("B76").Value = "eft:paltsek"
$newSheet.Range("C76").Value = "P8182"
$newSheet.Range("B77").Value = "eft:rinchen-tso"
$newSheet.Range("C77").Value = "P8273"
$newSheet.Range("B78").Value = "eft:manjusrigarbha"
$newSheet.Range("C78").Value = "P4CZ16780"
$newSheet.Range("B79").Value = "eft:siladharma"
$newSheet.Range("C79").Value = "https://lod.dila.edu.tw/resource.php?id=A000089"
$newSheet.Range("B80").Value = "eft:zhang-yesh-d-"
$newSheet.Range("C80").Value = "P8205"
$newSheet.Range("B81").Value = "eft:sherab-lekpa"
$newSheet.Range("C81").Value = "P4242"
$newSheet.Range("B82").Value = "eft:sakya-yesh-"
$newSheet.Range("C82").Value = "P3285"
$newSheet.Range("B83").Value = "eft:jinavara"
$newSheet.Range("C83").Value = "P0TMP098"
$newSheet.Range("B84").Value = "eft:trakpa-gyaltsen"
$newSheet.Range("C84").Value = "P2637"
$newSheet.Range("B85").Value = "eft:phakpa-sherab"
$newSheet.Range("C85").Value = "P3709"
$newSheet.Range("B86").Value = "eft:kumarakalasa"
$newSheet.Range("C86").Value = "P4CZ15137"
$newSheet.Range("B87").Value = "eft:dipamkara-srijnana"
$newSheet.Range("C87").Value = "P3379"
$newSheet.Range("B88").Value = "eft:pa-tshab-nyi-ma-grags"
$newSheet.Range("C88").Value = "P5651"
$newSheet.Range("B89").Value = "eft:band-yesh-de"
$newSheet.Range("C89").Value = "P8205"
$newSheet.Range("B90").Value = "eft:buddhakaravarma"
$newSheet.Range("C90").Value = "P8245"

# Select C1 on the new sheet (matches diff: <selection activeCell="C1" sqref="C1"/>)
$newSheet.Range("C1").Select()
